$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B142").Value = 63902
$ws.Range("E142").Value = 34.04
$ws.Range("F142").Value = 2
$ws.Range("G142").Value = 64.04000000000001
$ws.Range("B143").Value = 48654
$ws.Range("E143").Value = 38.26
$ws.Range("F143").Value = -1
$ws.Range("G143").Value = -32.02
$ws.Range("B154").Value = 64350
$ws.Range("E154").Value = 70.63
$ws.Range("F154").Value = 101
$ws.Range("G154").Value = 6710.44
$ws.Range("B155").Value = 57756
$ws.Range("E155").Value = 79.37
$ws.Range("F155").Value = -100
$ws.Range("G155").Value = -6644
$ws.Range("B156").Value = 53925
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 66.44
$ws.Range("B271").Value = 64973
$ws.Range("E271").Value = 35.4
$ws.Range("F271").Value = 150
$ws.Range("G271").Value = 4995
$ws.Range("B272").Value = 48706
$ws.Range("E272").Value = 39.8
$ws.Range("F272").Value = -144
$ws.Range("G272").Value = -4795.2
$ws.Range("B305").Value = 62997
$ws.Range("F305").Value = 72
$ws.Range("G305").Value = 22020.48
$ws.Range("B306").Value = 57854
$ws.Range("F306").Value = 2
$ws.Range("G306").Value = 611.6799999999999
$ws.Range("B309").Value = 61610
$ws.Range("E309").Value = 122.71
$ws.Range("F309").Value = -58
$ws.Range("G309").Value = -5957.18
$ws.Range("B310").Value = 63565
$ws.Range("E310").Value = 109.19
$ws.Range("F310").Value = 60
$ws.Range("G310").Value = 6162.6
$ws.Range("B342").Value = 57802
$ws.Range("E342").Value = 162.71
$ws.Range("F342").Value = -79
$ws.Range("G342").Value = -11334.92
$ws.Range("B344").Value = 63531
$ws.Range("E344").Value = 152.53
$ws.Range("F344").Value = 80
$ws.Range("G344").Value = 11478.4
$ws.Range("B347").Value = 55356
$ws.Range("E347").Value = 54.04
$ws.Range("F347").Value = -158
$ws.Range("G347").Value = -7527.12
$ws.Range("B348").Value = 63510
$ws.Range("E348").Value = 50.66
$ws.Range("F348").Value = 167
$ws.Range("G348").Value = 7955.88
$ws.Range("B371").Value = 63564
$ws.Range("E371").Value = 137.16
$ws.Range("F371").Value = 57
$ws.Range("G371").Value = 7353.57
$ws.Range("B372").Value = 61608
$ws.Range("E372").Value = 154.12
$ws.Range("F372").Value = -56
$ws.Range("G372").Value = -7224.56
$ws.Range("B374").Value = 63560
$ws.Range("E374").Value = 134.87
$ws.Range("F374").Value = 104
$ws.Range("G374").Value = 13193.44
$ws.Range("B375").Value = 60325
$ws.Range("E375").Value = 151.57
$ws.Range("F375").Value = -102
$ws.Range("G375").Value = -12939.72
$ws.Range("B381").Value = 57817
$ws.Range("F381").Value = 3
$ws.Range("G381").Value = 239.43
$ws.Range("B382").Value = 62865
$ws.Range("F382").Value = 151
$ws.Range("G382").Value = 12051.31
$ws.Range("B392").Value = 62933
$ws.Range("F392").Value = 146
$ws.Range("G392").Value = 8632.98
$ws.Range("B393").Value = 57835
$ws.Range("F393").Value = 1
$ws.Range("G393").Value = 59.13
$ws.Range("B411").Value = 57856
$ws.Range("F411").Value = 2
$ws.Range("G411").Value = 342.66
$ws.Range("B412").Value = 63007
$ws.Range("F412").Value = 984
$ws.Range("G412").Value = 168588.72
$ws.Range("B578").Value = 64915
$ws.Range("E578").Value = 20.98
$ws.Range("F578").Value = 40
$ws.Range("G578").Value = 789.2
$ws.Range("B579").Value = 45695
$ws.Range("E579").Value = 23.58
$ws.Range("F579").Value = -36
$ws.Range("G579").Value = -710.28
$ws.Range("B582").Value = 64922
$ws.Range("E582").Value = 20.98
$ws.Range("F582").Value = 207
$ws.Range("G582").Value = 4084.11
$ws.Range("B583").Value = 45706
$ws.Range("E583").Value = 23.58
$ws.Range("F583").Value = -202
$ws.Range("G583").Value = -3985.46
$ws.Range("B585").Value = 45718
$ws.Range("E585").Value = 19.38
$ws.Range("F585").Value = -294
$ws.Range("G585").Value = -4768.68
$ws.Range("B586").Value = 64927
$ws.Range("E586").Value = 17.26
$ws.Range("F586").Value = 295
$ws.Range("G586").Value = 4784.9
$ws.Range("B701").Value = 64833
$ws.Range("E701").Value = 34.9
$ws.Range("F701").Value = 99
$ws.Range("G701").Value = 3250.17
$ws.Range("B702").Value = 60025
$ws.Range("E702").Value = 37.22
$ws.Range("F702").Value = -98
$ws.Range("G702").Value = -3217.34
$ws.Range("B712").Value = 64830
$ws.Range("E712").Value = 34.9
$ws.Range("F712").Value = 117
$ws.Range("G712").Value = 3841.11
$ws.Range("B713").Value = 60022
$ws.Range("E713").Value = 37.22
$ws.Range("F713").Value = -113
$ws.Range("G713").Value = -3709.79
$ws.Range("B864").Value = 65079
$ws.Range("E864").Value = 43.44
$ws.Range("F864").Value = 21
$ws.Range("G864").Value = 858.27
$ws.Range("B865").Value = 54751
$ws.Range("E865").Value = 46.34
$ws.Range("F865").Value = -19
$ws.Range("G865").Value = -776.53
